$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q, continuing the
# 0..13 sequence already present in B1:O1.
# First copy the style of O1 (bold/border/centered) onto P1:Q1, then set values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Swap values between columns I/K and M/O for data rows 2-25
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# Fill the two new columns (P, Q) for rows 2-25 with value 2
$ws.Range("P2:Q25").Value = 2
